# Update "report co so" workbook:
#  1. CHI TIẾT DOANH THU: reorder columns A/B/C so that
#     A = Ngày thực hiện, B = Tiền tố, C = Mã dịch vụ (was A=Tiền tố, B=Mã dịch vụ, C=Ngày thực hiện)
#  2. Insert a brand-new sheet "CHI TIẾT CHI TIÊU" right after "CHI TIẾT DOANH THU"
#     with a detailed list of expense transactions.
#  3. Rename "DAONH SỐ CÁ NHÂN" -> "DOANH SỐ CÁ NHÂN" (typo fix).
#  4. Rename "CHI TIÊU" -> "CHI TIÊU TỔNG HỢP".
#  5. "LŨY KẾ NGÀY" is left untouched.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) CHI TIẾT DOANH THU - rotate columns A,B,C (header row + 8 data rows)
# ---------------------------------------------------------------------------
$wsDoanhThu = $wb.Worksheets.Item("CHI TIẾT DOANH THU")

for ($r = 1; $r -le 9; $r++) {
    $oldA = $wsDoanhThu.Cells.Item($r, 1).Value()
    $oldB = $wsDoanhThu.Cells.Item($r, 2).Value()
    $oldC = $wsDoanhThu.Cells.Item($r, 3).Value()

    # New A = old C (date text) - briefly force text format so the string
    # doesn't get auto-converted into a date serial, then drop back to the
    # default "Normal" style so no stray formatting is left behind.
    $cellA = $wsDoanhThu.Cells.Item($r, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $oldC
    $cellA.Style = "Normal"
    # New B = old A (Tiền tố)
    $wsDoanhThu.Cells.Item($r, 2).Value = $oldA
    # New C = old B (Mã dịch vụ)
    $wsDoanhThu.Cells.Item($r, 3).Value = $oldB
}

# ---------------------------------------------------------------------------
# 2) Insert new sheet "CHI TIẾT CHI TIÊU" right after "CHI TIẾT DOANH THU"
# ---------------------------------------------------------------------------
$wsChiTietChiTieu = $wb.Worksheets.Add($null, $wsDoanhThu)
$wsChiTietChiTieu.Name = "CHI TIẾT CHI TIÊU"

$headers = @("Tiền tố", "Mã chi tiêu", "Ngày chi", "Cơ sở", "Phân loại", "Lượng chi")
for ($c = 1; $c -le 6; $c++) {
    $wsChiTietChiTieu.Cells.Item(1, $c).Value = $headers[$c - 1]
}

$rows = @(
    @("CT", 584, "07-02-2024", "SÓC TRĂNG", "Chi Phí CTV", 5000000),
    @("CT", 585, "07-02-2024", "SÓC TRĂNG", "Chí Phí Bác Sĩ Ngoài", 6250000),
    @("CT", 586, "07-02-2024", "SÓC TRĂNG", "Phúc lợi công ty", 400000),
    @("CT", 587, "07-02-2024", "SÓC TRĂNG", "Chi Phí Sinh Hoạt Tại Cơ Sở", 300000),
    @("CT", 588, "07-02-2024", "SÓC TRĂNG", "Chi Phí Vận Hành", 1000000),
    @("CT", 589, "07-03-2024", "SÓC TRĂNG", "Chi Phí Sinh Hoạt Tại Cơ Sở", 1506000),
    @("CT", 600, "07-05-2024", "SÓC TRĂNG", "Chi Phí Vận Hành", 1300000),
    @("CT", 601, "07-05-2024", "SÓC TRĂNG", "Chi Phí CTV", 11000000),
    @("CT", 602, "07-05-2024", "SÓC TRĂNG", "Chi Phí CTV", 1650000),
    @("CT", 603, "07-05-2024", "SÓC TRĂNG", "Chi Phí Hạ Tầng", 4544000),
    @("CT", 604, "07-05-2024", "SÓC TRĂNG", "Phúc lợi công ty", 200000),
    @("CT", 605, "07-05-2024", "SÓC TRĂNG", "Trang thiết bị Y Tế", 3500000),
    @("CT", 606, "07-05-2024", "SÓC TRĂNG", "Chi Phí Sinh Hoạt Tại Cơ Sở", 782000),
    @("CT", 607, "07-07-2024", "SÓC TRĂNG", "Chi Phí Sinh Hoạt Tại Cơ Sở", 285000)
)

$rowIndex = 2
foreach ($row in $rows) {
    $wsChiTietChiTieu.Cells.Item($rowIndex, 1).Value = $row[0]
    $wsChiTietChiTieu.Cells.Item($rowIndex, 2).Value = $row[1]
    $cDate = $wsChiTietChiTieu.Cells.Item($rowIndex, 3)
    $cDate.NumberFormat = "@"
    $cDate.Value = $row[2]
    $cDate.Style = "Normal"
    $wsChiTietChiTieu.Cells.Item($rowIndex, 4).Value = $row[3]
    $wsChiTietChiTieu.Cells.Item($rowIndex, 5).Value = $row[4]
    $wsChiTietChiTieu.Cells.Item($rowIndex, 6).Value = $row[5]
    $rowIndex++
}

# ---------------------------------------------------------------------------
# 3) Rename "DAONH SỐ CÁ NHÂN" -> "DOANH SỐ CÁ NHÂN"
# ---------------------------------------------------------------------------
$wsCaNhan = $wb.Worksheets.Item("DAONH SỐ CÁ NHÂN")
$wsCaNhan.Name = "DOANH SỐ CÁ NHÂN"

# ---------------------------------------------------------------------------
# 4) Rename "CHI TIÊU" -> "CHI TIÊU TỔNG HỢP"
# ---------------------------------------------------------------------------
$wsChiTieu = $wb.Worksheets.Item("CHI TIÊU")
$wsChiTieu.Name = "CHI TIÊU TỔNG HỢP"
